# Update leaderboard standings on Sheet1 (totaalstand_EL1_EL8.xlsx)
# Applies the new match results: several players' rows are refreshed with
# updated stats (and a few players swap ranking positions / rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Yannick den Daggelder ---
$ws.Range("C2").Value = 64
$ws.Range("F2").Value = 47556
$ws.Range("G2").Value = 2710
$ws.Range("H2").Value = 52.65
$ws.Range("I2").Value = 69

# --- Row 3: now Rocky Van Den Eeckhoudt ---
$ws.Range("B3").Value = "Rocky Van Den Eeckhoudt"
$ws.Range("C3").Value = 35
$ws.Range("D3").Value = 5
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 20188
$ws.Range("G3").Value = 1102
$ws.Range("H3").Value = 54.96
$ws.Range("I3").Value = 40
$ws.Range("J3").Value = 1

# --- Row 4: now Dartin Dan ---
$ws.Range("B4").Value = "Dartin Dan"
$ws.Range("C4").Value = 31
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 19917
$ws.Range("G4").Value = 1150
$ws.Range("H4").Value = 51.96
$ws.Range("I4").Value = 34
$ws.Range("J4").Value = 0

# --- Row 8: Robin Willis ---
$ws.Range("C8").Value = 24
$ws.Range("F8").Value = 20654
$ws.Range("G8").Value = 1272
$ws.Range("H8").Value = 48.71
$ws.Range("I8").Value = 24

# --- Row 9: now Noah B ---
$ws.Range("B9").Value = "Noah B"
$ws.Range("C9").Value = 22
$ws.Range("F9").Value = 14992
$ws.Range("G9").Value = 862
$ws.Range("H9").Value = 52.18
$ws.Range("I9").Value = 22
$ws.Range("J9").Value = 1

# --- Row 10: now Milan Schoenmakers ---
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Milan Schoenmakers"
$ws.Range("C10").Value = 20
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 11242
$ws.Range("G10").Value = 584
$ws.Range("H10").Value = 57.75
$ws.Range("I10").Value = 21

# --- Row 11: now Max Walter ---
$ws.Range("B11").Value = "Max Walter"
$ws.Range("C11").Value = 18
$ws.Range("D11").Value = 0
$ws.Range("F11").Value = 9206
$ws.Range("G11").Value = 663
$ws.Range("H11").Value = 41.66
$ws.Range("I11").Value = 18

# --- Row 12: now Louis Tweddle ---
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Louis Tweddle"
$ws.Range("C12").Value = 17
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 11864
$ws.Range("G12").Value = 581
$ws.Range("H12").Value = 61.26
$ws.Range("I12").Value = 18

# --- Row 13: now Afendi Kelana ---
$ws.Range("B13").Value = "Afendi Kelana"
$ws.Range("C13").Value = 8
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 8344
$ws.Range("G13").Value = 432
$ws.Range("H13").Value = 57.94
$ws.Range("I13").Value = 11

# --- Row 14: now Diego Meerveld ---
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "Diego Meerveld"
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 0
$ws.Range("F14").Value = 7022
$ws.Range("G14").Value = 457
$ws.Range("H14").Value = 46.1

# --- Row 15: now Francesco Petruccelli ---
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Francesco Petruccelli"
$ws.Range("F15").Value = 6404
$ws.Range("G15").Value = 349
$ws.Range("H15").Value = 55.05

# --- Row 16: now martin Berry ---
$ws.Range("A16").Value = 13
$ws.Range("B16").Value = "martin Berry"
$ws.Range("C16").Value = 8
$ws.Range("D16").Value = 1
$ws.Range("F16").Value = 4295
$ws.Range("G16").Value = 285
$ws.Range("H16").Value = 45.21
$ws.Range("I16").Value = 9

# --- Row 17: now Magnus Gladh ---
$ws.Range("A17").Value = 16
$ws.Range("B17").Value = "Magnus Gladh"
$ws.Range("F17").Value = 6401
$ws.Range("G17").Value = 377
$ws.Range("H17").Value = 50.94
